$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the 2021 year column header (O4), copying the formatting from N4 (2020)
$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)
$ws.Range("O4").Value2 = 2021

# Add the corresponding data value (O5), copying the formatting from N5
$ws.Range("N5").Copy()
$ws.Range("O5").PasteSpecial(-4122)
$ws.Range("O5").Value2 = 1.5020015556876996

$excel.CutCopyMode = $false

# Update the selection to match the saved view state
$ws.Range("Q5").Select()
